$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.770.37'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.56%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.720.89'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.31%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.0000'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.31%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.10'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.45%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.26%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4762'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -2.20%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2552'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.95%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06110'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.94%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.716.71'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.15%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '15.84'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +2.45%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.06900'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.67%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.5949'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.35%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.405'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -1.52%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '76.34'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.11%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.001'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.39%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.678.77'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.27%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.000'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.36%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007015'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.97%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.24'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.02%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.939.93'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.09%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.362'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.09%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.296'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -1.56%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.052'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.06%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '140.63'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +3.28%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.11'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.22%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.778'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +2.92%  '
$ws.Range('B28').Value = 'BitcoinCash'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '105.79'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.27%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.369'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -2.04%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.939'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +1.86%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.07877'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -1.02%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.609'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.00%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04533'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +2.68%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.586'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.63%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9891'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.42%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6097'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -1.17%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9209'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -1.51%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.482'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +4.51%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.957'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -1.58%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9997'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.30%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.677'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +4.65%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.01477'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.41%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '100.09'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.50%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.3774'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.85%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '6.698'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -1.81%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1141'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.80%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05344'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.05%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.763'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.38%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '29.54'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -2.90%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.227'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +1.15%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '50.58'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.92%  '
